$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.443.89"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "2.310.30"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.42%  "
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("E11").Value = "  +3.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("D15").Value = "2.669.55"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.309.46"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.811"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").Value = "43.349.69"
$ws.Range("E19").Value = "  +3.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "0.0₃0930"
$ws.Range("E21").Value = "  +3.45%  "
$ws.Range("E22").Value = "  +3.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "242.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.76%  "
$ws.Range("E25").Value = "  +2.79%  "
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +6.04%  "
$ws.Range("E29").Value = "  +7.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.15%  "
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("E39").Value = "  +4.46%  "
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.66%  "
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0298"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.64%  "
$ws.Range("D45").Value = "1.987.29"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("E51").Value = "  +9.59%  "
